$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("g2.3")

# Update the "Quantidade 2024/2014" -> "Quantidade 2024/2015" and
# "Valor 2024/2014" -> "Valor 2024/2015" labels in column B, and the
# corresponding recalculated values in column D, for rows 2-13.
# (Min year used in the comparison changed from 2014 to 2015.)

$ws.Range("B2").Value = "Quantidade 2024/2015"
$ws.Range("D2").Value = 42.03766661270304

$ws.Range("B3").Value = "Quantidade 2024/2015"
$ws.Range("D3").Value = 84.08918303135511

$ws.Range("B4").Value = "Quantidade 2024/2015"
$ws.Range("D4").Value = -14.5014850570501

$ws.Range("B5").Value = "Quantidade 2024/2015"
$ws.Range("D5").Value = -55.51025561441625

$ws.Range("B6").Value = "Quantidade 2024/2015"
$ws.Range("D6").Value = 254.5842217484009

$ws.Range("B7").Value = "Quantidade 2024/2015"
$ws.Range("D7").Value = 301.0989010989011

$ws.Range("B8").Value = "Valor 2024/2015"
$ws.Range("D8").Value = 70.77527442490312

$ws.Range("B9").Value = "Valor 2024/2015"
$ws.Range("D9").Value = 44.24612788095496

$ws.Range("B10").Value = "Valor 2024/2015"
$ws.Range("D10").Value = -22.80201896534912

$ws.Range("B11").Value = "Valor 2024/2015"
$ws.Range("D11").Value = -52.00431130021104

$ws.Range("B12").Value = "Valor 2024/2015"
$ws.Range("D12").Value = 311.1692808246118

$ws.Range("B13").Value = "Valor 2024/2015"
$ws.Range("D13").Value = 96.32603784185962

$wb.Save()
